# Update "想去人数" (number of people interested) counts across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 261
$ws1.Range("F5").Value = 3082
$ws1.Range("F6").Value = 2056
$ws1.Range("F9").Value = 1152
$ws1.Range("F11").Value = 897
$ws1.Range("F12").Value = 74

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 29

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 261
$ws4.Range("F5").Value = 3082
$ws4.Range("F6").Value = 2056
$ws4.Range("F8").Value = 29
$ws4.Range("F10").Value = 1152
$ws4.Range("F12").Value = 897
$ws4.Range("F13").Value = 74
